$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 9533944
$ws.Range("J17").Value = 10267132
$ws.Range("L17").Value = 30801396
$ws.Range("N17").Value = -30801732
$ws.Range("H32").Value = 4000
$ws.Range("J32").Value = 4000
$ws.Range("L32").Value = 4000
$ws.Range("N32").Value = -4652
$ws.Range("H33").Value = 1421.4242
$ws.Range("I33").Value = 138.90909
$ws.Range("J33").Value = 2062.682
$ws.Range("K33").Value = 138.90909
$ws.Range("L33").Value = 2062.682
$ws.Range("M33").Value = 90.09091000000001
$ws.Range("N33").Value = -2520.682
$ws.Range("H40").Value = 15792420
$ws.Range("I40").Value = 3368.3076
$ws.Range("J40").Value = 50002030
$ws.Range("K40").Value = 3368.3076
$ws.Range("L40").Value = 50002030
$ws.Range("M40").Value = -3193.3076
$ws.Range("N40").Value = -50002380
$ws.Range("H51").Value = 20783.357
$ws.Range("I51").Value = 21193
$ws.Range("K51").Value = 21193
$ws.Range("M51").Value = -20709
$ws.Range("H86").Value = 3819.1428
$ws.Range("J86").Value = 3554.1428
$ws.Range("L86").Value = 3554.1428
$ws.Range("N86").Value = -5800.1428
$ws.Range("H89").Value = 3819.1428
$ws.Range("J89").Value = 3554.1428
$ws.Range("L89").Value = 17770.714
$ws.Range("N89").Value = -29002.714
$ws.Range("H112").Value = 6971696
$ws.Range("J112").Value = 6971696
$ws.Range("L112").Value = 20915088
$ws.Range("N112").Value = -20917304
$ws.Range("H135").Value = 1714.3871
$ws.Range("I135").Value = 1368.8462
$ws.Range("K135").Value = 12319.6158
$ws.Range("M135").Value = -9784.6158
$ws.Range("H137").Value = 23969.09
$ws.Range("I137").Value = 38096.613
$ws.Range("K137").Value = 114289.839
$ws.Range("M137").Value = -111739.839
$ws.Range("H141").Value = 2591.625
$ws.Range("I141").Value = 2274.2307
$ws.Range("K141").Value = 6822.6921
$ws.Range("M141").Value = -1642.6921

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H14").Value = 6001
$ws.Range("J14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("N14").ClearContents()
$ws.Range("H30").Value = 8736.333000000001
$ws.Range("I30").Value = 12504.5
$ws.Range("J30").Value = 1200
$ws.Range("K30").Value = 12504.5
$ws.Range("L30").Value = 1200
$ws.Range("M30").Value = -12354.5
$ws.Range("N30").Value = -1500
$ws.Range("H32").Value = 4862.769
$ws.Range("I32").Value = 2909.9265
$ws.Range("K32").Value = 2909.9265
$ws.Range("M32").Value = -2622.9265
$ws.Range("H74").Value = 84876.414
$ws.Range("I74").Value = 125158.375
$ws.Range("K74").Value = 125158.375
$ws.Range("M74").Value = -124284.375
$ws.Range("H77").Value = 84876.414
$ws.Range("I77").Value = 125158.375
$ws.Range("K77").Value = 625791.875
$ws.Range("M77").Value = -621423.875
$ws.Range("H110").Value = 13554.272
$ws.Range("I110").Value = 14074.625
$ws.Range("J110").Value = 12166.667
$ws.Range("K110").Value = 14074.625
$ws.Range("L110").Value = 12166.667
$ws.Range("M110").Value = -12029.625
$ws.Range("N110").Value = -16256.667
$ws.Range("H132").Value = 2996.3333
$ws.Range("I132").Value = 2996.3333
$ws.Range("K132").Value = 8988.999899999999
$ws.Range("M132").Value = -6458.999899999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H138").Value = 69967
$ws.Range("J138").Value = 69967
$ws.Range("L138").Value = 69967
$ws.Range("N138").Value = -80247

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 436361.3
$ws.Range("I31").Value = 715202
$ws.Range("J31").Value = 2609.111
$ws.Range("K31").Value = 715202
$ws.Range("L31").Value = 2609.111
$ws.Range("M31").Value = -714907
$ws.Range("N31").Value = -3199.111
$ws.Range("H34").Value = 436361.3
$ws.Range("I34").Value = 715202
$ws.Range("J34").Value = 2609.111
$ws.Range("K34").Value = 715202
$ws.Range("L34").Value = 2609.111
$ws.Range("M34").Value = -715000
$ws.Range("N34").Value = -3013.111
$ws.Range("H58").Value = 2318.9
$ws.Range("I58").Value = 2574.4
$ws.Range("J58").Value = 2063.4
$ws.Range("K58").Value = 2574.4
$ws.Range("L58").Value = 2063.4
$ws.Range("M58").Value = -2371.4
$ws.Range("N58").Value = -2469.4
$ws.Range("H64").Value = 35125.184
$ws.Range("J64").Value = 35125.184
$ws.Range("L64").Value = 35125.184
$ws.Range("N64").Value = -35621.184
$ws.Range("H67").Value = 35125.184
$ws.Range("J67").Value = 35125.184
$ws.Range("L67").Value = 35125.184
$ws.Range("N67").Value = -36841.184
$ws.Range("H68").Value = 40071.25
$ws.Range("J68").Value = 40071.25
$ws.Range("L68").Value = 40071.25
$ws.Range("N68").Value = -41569.25
$ws.Range("H71").Value = 40071.25
$ws.Range("J71").Value = 40071.25
$ws.Range("L71").Value = 120213.75
$ws.Range("N71").Value = -127701.75
$ws.Range("H136").Value = 2318.9
$ws.Range("I136").Value = 2574.4
$ws.Range("J136").Value = 2063.4
$ws.Range("K136").Value = 7723.200000000001
$ws.Range("L136").Value = 6190.200000000001
$ws.Range("M136").Value = -5173.200000000001
$ws.Range("N136").Value = -11290.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 834.9091
$ws.Range("J23").Value = 943.8946999999999
$ws.Range("L23").Value = 2831.6841
$ws.Range("N23").Value = -3301.6841
$ws.Range("H109").Value = 143890
$ws.Range("I109").Value = 167455
$ws.Range("J109").Value = 2500
$ws.Range("K109").Value = 502365
$ws.Range("L109").Value = 7500
$ws.Range("M109").Value = -501325
$ws.Range("N109").Value = -9580
$ws.Range("H121").Value = 125004824
$ws.Range("I121").Value = 500003500
$ws.Range("J121").Value = 5266
$ws.Range("K121").Value = 1500010500
$ws.Range("L121").Value = 15798
$ws.Range("M121").Value = -1500009190
$ws.Range("N121").Value = -18418

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H104").Value = 74053.8
$ws.Range("J104").Value = 74053.8
$ws.Range("L104").Value = 74053.8
$ws.Range("N104").Value = -81041.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1298.5526
$ws.Range("I22").Value = 1038.1578
$ws.Range("J22").Value = 1558.9474
$ws.Range("K22").Value = 1038.1578
$ws.Range("L22").Value = 1558.9474
$ws.Range("M22").Value = -743.1578
$ws.Range("N22").Value = -2148.9474
$ws.Range("H27").Value = 1298.5526
$ws.Range("I27").Value = 1038.1578
$ws.Range("J27").Value = 1558.9474
$ws.Range("K27").Value = 1038.1578
$ws.Range("L27").Value = 1558.9474
$ws.Range("M27").Value = -931.1578
$ws.Range("N27").Value = -1772.9474
$ws.Range("H93").Value = 3693.1667
$ws.Range("I93").Value = 2303.7778
$ws.Range("K93").Value = 2303.7778
$ws.Range("M93").Value = -1055.7778

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H80").Value = 14400
$ws.Range("J80").Value = 14400
$ws.Range("L80").Value = 14400
$ws.Range("N80").Value = -16396
$ws.Range("H83").Value = 14400
$ws.Range("J83").Value = 14400
$ws.Range("L83").Value = 43200
$ws.Range("N83").Value = -53184
$ws.Range("H96").Value = 5357.3335
$ws.Range("J96").Value = 6082.4
$ws.Range("L96").Value = 6082.4
$ws.Range("N96").Value = -8828.4
